# break out stock.yaml completed
# - convert the bsecode column (D) on rows 222:227 of the "day" sheet from
#   text to real numbers
# - append 7 new rows (228:234) of freshly scraped stock data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. D222:D227 were stored as text ("532777" ...); make them numeric ---
$bsecodes = @{
    222 = 532777
    223 = 517354
    224 = 532174
    225 = 543066
    226 = 532898
    227 = 532720
}
foreach ($r in $bsecodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $bsecodes[$r]
}

# --- 2. append the new rows scraped on 29/07/2024 ---
# columns: sr, nsecode, name, bsecode, per_chg, close, volume, timeframe, Date Time
$newRows = @(
    @(1, "PERSISTENT", "Persistent Systems Limited", "533179", 0.12, 4778.6, 438646, "day", "29/07/2024 11:35:45"),
    @(2, "ASIANPAINT", "Asian Paints Limited", "500820", 0.15, 2954.7, 854461, "day", "29/07/2024 11:35:45"),
    @(3, "M&M", "Mahindra & Mahindra Limited", "500520", 1.57, 2933, 2388937, "day", "29/07/2024 11:35:45"),
    @(4, "AXISBANK", "Axis Bank Limited", "532215", -0.62, 1170.05, 19598573, "day", "29/07/2024 11:35:45"),
    @(5, "IRCTC", "Indian Railway Catering & Tourism Corporation Ltd", "542830", 0.74, 991.45, 1762580, "day", "29/07/2024 11:35:45"),
    @(6, "RECLTD", "Rural Electrification Corporation Limited", "532955", 2.99, 644.6, 23006702, "day", "29/07/2024 11:35:45"),
    @(7, "PFC", "Power Finance Corporation Limited", "532810", 2.58, 552.85, 11024457, "day", "29/07/2024 11:35:45")
)

$startRow = 228
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    # bsecode keeps coming through as text - force it to stay text by typing
    # it the way a user would (leading apostrophe), same as the rest of the
    # sheet's bsecode column before this commit.
    $ws.Cells.Item($row, 4).Value = "'" + $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 7).Value = $data[6]
    $ws.Cells.Item($row, 8).Value = $data[7]
    $ws.Cells.Item($row, 9).Value = $data[8]
}
